$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2025-11-02 Sunday"

# Update each math-expression cell in the table (row-major order matches
# document order: Cell(1,1), Cell(1,2), ... Cell(1,5), Cell(2,1), ...)
$t = $d.Tables.Item(1)
$values = @(
    "38-17=",
    "59-50=",
    "38-33=",
    "36-27=",
    "57-20=",
    "87-78=",
    "67+28=",
    "63-22=",
    "97-73=",
    "49-41=",
    "42+15=",
    "75-34=",
    "26+51=",
    "81-2=",
    "28+5=",
    "0+5=",
    "44-31=",
    "14+67=",
    "52+37=",
    "42+23=",
    "50-45=",
    "10-1=",
    "61-56=",
    "27+41=",
    "85-18=",
    "80+1=",
    "29+34=",
    "36+49=",
    "85-73=",
    "5+21=",
    "97-83=",
    "71-11=",
    "5+29=",
    "24-10=",
    "74+16=",
    "1+3=",
    "64-29=",
    "44-7=",
    "82-8=",
    "28+8=",
    "28+59=",
    "44+0=",
    "73+16=",
    "39+41=",
    "91-26=",
    "4+57=",
    "0+31=",
    "63-10=",
    "98-85=",
    "20+66=",
    "21-9=",
    "11+8=",
    "55+3=",
    "57+42=",
    "65+7=",
    "7+73=",
    "34-24=",
    "27-16=",
    "40-7=",
    "56+37=",
    "88+10=",
    "61-34=",
    "29+68=",
    "15+57=",
    "61+33=",
    "84+3=",
    "93-81=",
    "4+81=",
    "79-69=",
    "27-24=",
    "56-11=",
    "10+12=",
    "11+54=",
    "70-47=",
    "77-53=",
    "87-45=",
    "3+81=",
    "53+12=",
    "34+52=",
    "35+13=",
    "84-79=",
    "87-86=",
    "63-62=",
    "1+87=",
    "93-60=",
    "34-17=",
    "22+48=",
    "76-41=",
    "96-41=",
    "59-28=",
    "74+1=",
    "72-59=",
    "29+59=",
    "67-43=",
    "78-68=",
    "7+37=",
    "51-50=",
    "24+14=",
    "87-70=",
    "86-57="
)

$cols = 5
$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Updated" $idx "table cells and the date heading."
